# Auto-generated edit script applying numeric corrections to profit-calculation
# columns (H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across several
# sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 118.55556
$ws.Cells.Item(33, 9).Value = 115.6
$ws.Cells.Item(33, 11).Value = 115.6
$ws.Cells.Item(33, 13).Value = 113.4
$ws.Cells.Item(64, 8).Value = 3718.3274
$ws.Cells.Item(64, 9).Value = 4000
$ws.Cells.Item(64, 10).Value = 3702.077
$ws.Cells.Item(64, 11).Value = 4000
$ws.Cells.Item(64, 12).Value = 3702.077
$ws.Cells.Item(64, 13).Value = -3752
$ws.Cells.Item(64, 14).Value = -4198.077
$ws.Cells.Item(67, 8).Value = 3718.3274
$ws.Cells.Item(67, 9).Value = 4000
$ws.Cells.Item(67, 10).Value = 3702.077
$ws.Cells.Item(67, 11).Value = 4000
$ws.Cells.Item(67, 12).Value = 3702.077
$ws.Cells.Item(67, 13).Value = -3142
$ws.Cells.Item(67, 14).Value = -5418.077
$ws.Cells.Item(74, 8).Value = 5245.5713
$ws.Cells.Item(74, 9).Value = 5221.636
$ws.Cells.Item(74, 10).Value = 5333.3335
$ws.Cells.Item(74, 11).Value = 5221.636
$ws.Cells.Item(74, 12).Value = 5333.3335
$ws.Cells.Item(74, 13).Value = -4285.636
$ws.Cells.Item(74, 14).Value = -7205.3335
$ws.Cells.Item(77, 8).Value = 5245.5713
$ws.Cells.Item(77, 9).Value = 5221.636
$ws.Cells.Item(77, 10).Value = 5333.3335
$ws.Cells.Item(77, 11).Value = 26108.18
$ws.Cells.Item(77, 12).Value = 26666.6675
$ws.Cells.Item(77, 13).Value = -21428.18
$ws.Cells.Item(77, 14).Value = -36026.6675
$ws.Cells.Item(113, 8).Value = 3155.04
$ws.Cells.Item(113, 9).Value = 1756.875
$ws.Cells.Item(113, 10).Value = 3813
$ws.Cells.Item(113, 11).Value = 1756.875
$ws.Cells.Item(113, 12).Value = 3813
$ws.Cells.Item(113, 13).Value = 1497.125
$ws.Cells.Item(113, 14).Value = -10321
$ws.Cells.Item(115, 8).Value = 2647.3684
$ws.Cells.Item(115, 9).Value = 284
$ws.Cells.Item(115, 10).Value = 3491.4285
$ws.Cells.Item(115, 11).Value = 852
$ws.Cells.Item(115, 12).Value = 10474.2855
$ws.Cells.Item(115, 13).Value = 715
$ws.Cells.Item(115, 14).Value = -13608.2855
$ws.Cells.Item(121, 8).Value = 1205
$ws.Cells.Item(121, 10).Value = 1189.4736
$ws.Cells.Item(121, 12).Value = 3568.4208
$ws.Cells.Item(121, 14).Value = -7062.4208
$ws.Cells.Item(138, 8).Value = 2028.9138
$ws.Cells.Item(138, 9).Value = 1295.6957
$ws.Cells.Item(138, 10).Value = 2510.743
$ws.Cells.Item(138, 11).Value = 3887.0871
$ws.Cells.Item(138, 12).Value = 7532.228999999999
$ws.Cells.Item(138, 13).Value = 1252.9129
$ws.Cells.Item(138, 14).Value = -17812.229

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 383.21155
$ws.Cells.Item(2, 9).Value = 366.5909
$ws.Cells.Item(2, 10).Value = 474.625
$ws.Cells.Item(2, 11).Value = 366.5909
$ws.Cells.Item(2, 12).Value = 474.625
$ws.Cells.Item(2, 13).Value = -253.5909
$ws.Cells.Item(2, 14).Value = -700.625
$ws.Cells.Item(32, 8).Value = 22991.178
$ws.Cells.Item(32, 9).Value = 8861.339
$ws.Cells.Item(32, 10).Value = 62811.637
$ws.Cells.Item(32, 11).Value = 8861.339
$ws.Cells.Item(32, 12).Value = 62811.637
$ws.Cells.Item(32, 13).Value = -8574.339
$ws.Cells.Item(32, 14).Value = -63385.637
$ws.Cells.Item(61, 8).Value = 2888.7292
$ws.Cells.Item(61, 9).Value = 2754.2
$ws.Cells.Item(61, 10).Value = 3561.375
$ws.Cells.Item(61, 11).Value = 2754.2
$ws.Cells.Item(61, 12).Value = 3561.375
$ws.Cells.Item(61, 13).Value = -2542.2
$ws.Cells.Item(61, 14).Value = -3985.375
$ws.Cells.Item(63, 8).Value = 1043609.2
$ws.Cells.Item(63, 9).Value = 1236151.8
$ws.Cells.Item(63, 10).Value = 3879.2
$ws.Cells.Item(63, 11).Value = 1236151.8
$ws.Cells.Item(63, 12).Value = 3879.2
$ws.Cells.Item(63, 13).Value = -1235465.8
$ws.Cells.Item(63, 14).Value = -5251.2
$ws.Cells.Item(66, 8).Value = 1043609.2
$ws.Cells.Item(66, 9).Value = 1236151.8
$ws.Cells.Item(66, 10).Value = 3879.2
$ws.Cells.Item(66, 11).Value = 6180759
$ws.Cells.Item(66, 12).Value = 19396
$ws.Cells.Item(66, 13).Value = -6177327
$ws.Cells.Item(66, 14).Value = -26260
$ws.Cells.Item(74, 8).Value = 945.2414
$ws.Cells.Item(74, 9).Value = 576.6
$ws.Cells.Item(74, 10).Value = 1764.4445
$ws.Cells.Item(74, 11).Value = 576.6
$ws.Cells.Item(74, 12).Value = 1764.4445
$ws.Cells.Item(74, 13).Value = 297.4
$ws.Cells.Item(74, 14).Value = -3512.4445
$ws.Cells.Item(77, 8).Value = 945.2414
$ws.Cells.Item(77, 9).Value = 576.6
$ws.Cells.Item(77, 10).Value = 1764.4445
$ws.Cells.Item(77, 11).Value = 2883
$ws.Cells.Item(77, 12).Value = 8822.2225
$ws.Cells.Item(77, 13).Value = 1485
$ws.Cells.Item(77, 14).Value = -17558.2225
$ws.Cells.Item(88, 8).Value = 8544.048000000001
$ws.Cells.Item(88, 9).Value = 1580.2
$ws.Cells.Item(88, 10).Value = 10720.25
$ws.Cells.Item(88, 11).Value = 1580.2
$ws.Cells.Item(88, 12).Value = 10720.25
$ws.Cells.Item(88, 13).Value = -1174.2
$ws.Cells.Item(88, 14).Value = -11532.25
$ws.Cells.Item(91, 8).Value = 8544.048000000001
$ws.Cells.Item(91, 9).Value = 1580.2
$ws.Cells.Item(91, 10).Value = 10720.25
$ws.Cells.Item(91, 11).Value = 1580.2
$ws.Cells.Item(91, 12).Value = 10720.25
$ws.Cells.Item(91, 13).Value = -176.2
$ws.Cells.Item(91, 14).Value = -13528.25
$ws.Cells.Item(110, 8).Value = 802.5263
$ws.Cells.Item(110, 9).Value = 783.0769
$ws.Cells.Item(110, 10).Value = 844.6667
$ws.Cells.Item(110, 11).Value = 783.0769
$ws.Cells.Item(110, 12).Value = 844.6667
$ws.Cells.Item(110, 13).Value = 1261.9231
$ws.Cells.Item(110, 14).Value = -4934.6667
$ws.Cells.Item(116, 8).Value = 383.21155
$ws.Cells.Item(116, 9).Value = 366.5909
$ws.Cells.Item(116, 10).Value = 474.625
$ws.Cells.Item(116, 11).Value = 366.5909
$ws.Cells.Item(116, 12).Value = 474.625
$ws.Cells.Item(116, 13).Value = 1927.4091
$ws.Cells.Item(116, 14).Value = -5062.625
$ws.Cells.Item(132, 8).Value = 2904.6155
$ws.Cells.Item(132, 9).Value = 1857.7858
$ws.Cells.Item(132, 10).Value = 4125.9165
$ws.Cells.Item(132, 11).Value = 5573.357400000001
$ws.Cells.Item(132, 12).Value = 12377.7495
$ws.Cells.Item(132, 13).Value = -3043.357400000001
$ws.Cells.Item(132, 14).Value = -17437.7495
$ws.Cells.Item(136, 8).Value = 2888.7292
$ws.Cells.Item(136, 9).Value = 2754.2
$ws.Cells.Item(136, 10).Value = 3561.375
$ws.Cells.Item(136, 11).Value = 8262.599999999999
$ws.Cells.Item(136, 12).Value = 10684.125
$ws.Cells.Item(136, 13).Value = -5712.599999999999
$ws.Cells.Item(136, 14).Value = -15784.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 383.21155
$ws.Cells.Item(3, 9).Value = 366.5909
$ws.Cells.Item(3, 10).Value = 474.625
$ws.Cells.Item(3, 11).Value = 366.5909
$ws.Cells.Item(3, 12).Value = 474.625
$ws.Cells.Item(3, 13).Value = -252.5909
$ws.Cells.Item(3, 14).Value = -702.625
$ws.Cells.Item(22, 8).Value = 814.6875
$ws.Cells.Item(22, 9).Value = 771.9231
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 771.9231
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -598.9231
$ws.Cells.Item(22, 14).Value = -1346
$ws.Cells.Item(80, 8).Value = 48530.12
$ws.Cells.Item(80, 9).Value = 100822.63
$ws.Cells.Item(80, 10).Value = 7443.143
$ws.Cells.Item(80, 11).Value = 100822.63
$ws.Cells.Item(80, 12).Value = 7443.143
$ws.Cells.Item(80, 13).Value = -99824.63
$ws.Cells.Item(80, 14).Value = -9439.143
$ws.Cells.Item(83, 8).Value = 48530.12
$ws.Cells.Item(83, 9).Value = 100822.63
$ws.Cells.Item(83, 10).Value = 7443.143
$ws.Cells.Item(83, 11).Value = 504113.15
$ws.Cells.Item(83, 12).Value = 37215.715
$ws.Cells.Item(83, 13).Value = -499121.15
$ws.Cells.Item(83, 14).Value = -47199.715
$ws.Cells.Item(134, 8).Value = 12959.74
$ws.Cells.Item(134, 9).Value = 17452.879
$ws.Cells.Item(134, 10).Value = 4237.7646
$ws.Cells.Item(134, 11).Value = 52358.637
$ws.Cells.Item(134, 12).Value = 12713.2938
$ws.Cells.Item(134, 13).Value = -49823.637
$ws.Cells.Item(134, 14).Value = -17783.2938

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5862.773
$ws.Cells.Item(31, 9).Value = 984.63635
$ws.Cells.Item(31, 10).Value = 30253.455
$ws.Cells.Item(31, 11).Value = 984.63635
$ws.Cells.Item(31, 12).Value = 30253.455
$ws.Cells.Item(31, 13).Value = -689.63635
$ws.Cells.Item(31, 14).Value = -30843.455
$ws.Cells.Item(34, 8).Value = 5862.773
$ws.Cells.Item(34, 9).Value = 984.63635
$ws.Cells.Item(34, 10).Value = 30253.455
$ws.Cells.Item(34, 11).Value = 984.63635
$ws.Cells.Item(34, 12).Value = 30253.455
$ws.Cells.Item(34, 13).Value = -782.63635
$ws.Cells.Item(34, 14).Value = -30657.455
$ws.Cells.Item(58, 8).Value = 14037.471
$ws.Cells.Item(58, 9).Value = 22241.766
$ws.Cells.Item(58, 10).Value = 3890.0527
$ws.Cells.Item(58, 11).Value = 22241.766
$ws.Cells.Item(58, 12).Value = 3890.0527
$ws.Cells.Item(58, 13).Value = -22038.766
$ws.Cells.Item(58, 14).Value = -4296.0527
$ws.Cells.Item(107, 8).Value = 526.1429000000001
$ws.Cells.Item(107, 9).Value = 487.77274
$ws.Cells.Item(107, 10).Value = 666.8333
$ws.Cells.Item(107, 11).Value = 487.77274
$ws.Cells.Item(107, 12).Value = 666.8333
$ws.Cells.Item(107, 13).Value = 1432.22726
$ws.Cells.Item(107, 14).Value = -4506.8333
$ws.Cells.Item(132, 8).Value = 5621.4443
$ws.Cells.Item(132, 9).Value = 1000.7857
$ws.Cells.Item(132, 10).Value = 10597.538
$ws.Cells.Item(132, 11).Value = 3002.3571
$ws.Cells.Item(132, 12).Value = 31792.614
$ws.Cells.Item(132, 13).Value = -472.3571000000002
$ws.Cells.Item(132, 14).Value = -36852.614
$ws.Cells.Item(134, 8).Value = 1454827.6
$ws.Cells.Item(134, 9).Value = 1037.2413
$ws.Cells.Item(134, 10).Value = 4466250.5
$ws.Cells.Item(134, 11).Value = 3111.7239
$ws.Cells.Item(134, 12).Value = 13398751.5
$ws.Cells.Item(134, 13).Value = -576.7239
$ws.Cells.Item(134, 14).Value = -13403821.5
$ws.Cells.Item(136, 8).Value = 14037.471
$ws.Cells.Item(136, 9).Value = 22241.766
$ws.Cells.Item(136, 10).Value = 3890.0527
$ws.Cells.Item(136, 11).Value = 66725.298
$ws.Cells.Item(136, 12).Value = 11670.1581
$ws.Cells.Item(136, 13).Value = -64175.298
$ws.Cells.Item(136, 14).Value = -16770.1581

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 32446.6
$ws.Cells.Item(70, 9).Value = 48955.547
$ws.Cells.Item(70, 10).Value = 4508.385
$ws.Cells.Item(70, 11).Value = 48955.547
$ws.Cells.Item(70, 12).Value = 4508.385
$ws.Cells.Item(70, 13).Value = -48685.547
$ws.Cells.Item(70, 14).Value = -5048.385
$ws.Cells.Item(73, 8).Value = 32446.6
$ws.Cells.Item(73, 9).Value = 48955.547
$ws.Cells.Item(73, 10).Value = 4508.385
$ws.Cells.Item(73, 11).Value = 48955.547
$ws.Cells.Item(73, 12).Value = 4508.385
$ws.Cells.Item(73, 13).Value = -48019.547
$ws.Cells.Item(73, 14).Value = -6380.385
$ws.Cells.Item(132, 8).Value = 23893.123
$ws.Cells.Item(132, 9).Value = 30295.115
$ws.Cells.Item(132, 10).Value = 7888.143
$ws.Cells.Item(132, 11).Value = 90885.345
$ws.Cells.Item(132, 12).Value = 23664.429
$ws.Cells.Item(132, 13).Value = -88355.345
$ws.Cells.Item(132, 14).Value = -28724.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 7533.3335
$ws.Cells.Item(82, 9).Value = 1500
$ws.Cells.Item(82, 10).Value = 8287.5
$ws.Cells.Item(82, 11).Value = 1500
$ws.Cells.Item(82, 12).Value = 8287.5
$ws.Cells.Item(82, 13).Value = -1139
$ws.Cells.Item(82, 14).Value = -9009.5
$ws.Cells.Item(85, 8).Value = 7533.3335
$ws.Cells.Item(85, 9).Value = 1500
$ws.Cells.Item(85, 10).Value = 8287.5
$ws.Cells.Item(85, 11).Value = 1500
$ws.Cells.Item(85, 12).Value = 8287.5
$ws.Cells.Item(85, 13).Value = -252
$ws.Cells.Item(85, 14).Value = -10783.5
$ws.Cells.Item(93, 8).Value = 1917.1904
$ws.Cells.Item(93, 9).Value = 2334.75
$ws.Cells.Item(93, 10).Value = 1360.4445
$ws.Cells.Item(93, 11).Value = 2334.75
$ws.Cells.Item(93, 12).Value = 1360.4445
$ws.Cells.Item(93, 13).Value = -1086.75
$ws.Cells.Item(93, 14).Value = -3856.4445
$ws.Cells.Item(100, 8).Value = 2537
$ws.Cells.Item(100, 9).Value = 1922.5555
$ws.Cells.Item(100, 11).Value = 1922.5555
$ws.Cells.Item(100, 13).Value = -1381.5555

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 727.375
$ws.Cells.Item(113, 9).Value = 585.8
$ws.Cells.Item(113, 10).Value = 1435.25
$ws.Cells.Item(113, 11).Value = 1757.4
$ws.Cells.Item(113, 12).Value = 4305.75
$ws.Cells.Item(113, 13).Value = 412.6000000000001
$ws.Cells.Item(113, 14).Value = -8645.75
$ws.Cells.Item(136, 8).Value = 18963658
$ws.Cells.Item(136, 9).Value = 10786519
$ws.Cells.Item(136, 10).Value = 38462990
$ws.Cells.Item(136, 11).Value = 32359557
$ws.Cells.Item(136, 12).Value = 115388970
$ws.Cells.Item(136, 13).Value = -32357007
$ws.Cells.Item(136, 14).Value = -115394070
